$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Tnf", "Tnfrsf1a", "ECs", 3, 1, 5.429001, 16.287003, 0.04230716253661782, 0.04239440107683373, 3, 1, 36.95112266666666, 110.853368, 0.1740115908809209, 0.1775751473829744, 200.607681908456, 1805.469137176104, 0.007361936658654564, 0.007528192019431678),
  @("ECs", "Tnf", "Tnfrsf1a", "FAPs", 3, 1, 5.429001, 16.287003, 0.04230716253661782, 0.04239440107683373, 3, 1, 51.86084433333334, 155.582533, 0.2442250025331967, 0.2492264486514428, 281.552575746511, 2533.973181718599, 0.01033246687767785, 0.01056580602308417),
  @("ECs", "Tnf", "Tnfrsf1a", "Inflammatory-Mac", 3, 1, 5.429001, 16.287003, 0.04230716253661782, 0.04239440107683373, 3, 1, 64.73785366666665, 194.213561, 0.3048658902295037, 0.3111091917238571, 351.4618722941869, 3163.156850647682, 0.0128980107698103, 0.01318928785263076),
  @("ECs", "Tnf", "Tnfrsf1a", "Resolving-Mac", 3, 1, 5.429001, 16.287003, 0.04230716253661782, 0.04239440107683373, 2, 1, 12.7841595, 25.568319, 0.06020363583370166, 0.04095769119761797, 69.4052147096595, 416.431288257957, 0.002547045006511765, 0.001736376787812919),
  @("ECs", "Tnf", "Tnfrsf1a", "MuSCs", 3, 1, 5.429001, 16.287003, 0.04230716253661782, 0.04239440107683373, 3, 1, 46.01464833333333, 138.043945, 0.216693880522677, 0.2211315210441077, 249.813571816315, 2248.322146346835, 0.00916770322396334, 0.009374738393874201),
  @("FAPs", "Tnf", "Tnfrsf1a", "ECs", 3, 1, 1.36117, 4.08351, 0.01060733649339319, 0.01062920911485442, 3, 1, 36.95112266666666, 110.853368, 0.1740115908809209, 0.1775751473829744, 50.29675964018667, 452.67083676168, 0.001845799498224597, 0.001887483375134729),
  @("FAPs", "Tnf", "Tnfrsf1a", "FAPs", 3, 1, 1.36117, 4.08351, 0.01060733649339319, 0.01062920911485442, 3, 1, 51.86084433333334, 155.582533, 0.2442250025331967, 0.2492264486514428, 70.59142548120336, 635.3228293308301, 0.00259057678196942, 0.002649080039668713),
  @("FAPs", "Tnf", "Tnfrsf1a", "Inflammatory-Mac", 3, 1, 1.36117, 4.08351, 0.01060733649339319, 0.01062920911485442, 3, 1, 64.73785366666665, 194.213561, 0.3048658902295037, 0.3111091917238571, 88.11922427545666, 793.0730184791099, 0.003233815083022216, 0.003306844656386214),
  @("FAPs", "Tnf", "Tnfrsf1a", "Resolving-Mac", 3, 1, 1.36117, 4.08351, 0.01060733649339319, 0.01062920911485442, 2, 1, 12.7841595, 25.568319, 0.06020363583370166, 0.04095769119761797, 17.40141438661501, 104.40848631969, 0.0006386002234137773, 0.0004353478646011137),
  @("FAPs", "Tnf", "Tnfrsf1a", "MuSCs", 3, 1, 1.36117, 4.08351, 0.01060733649339319, 0.01062920911485442, 3, 1, 46.01464833333333, 138.043945, 0.216693880522677, 0.2211315210441077, 62.63375887188334, 563.70382984695, 0.002298544906763175, 0.002350453179063653),
  @("Inflammatory-Mac", "Tnf", "Tnfrsf1a", "ECs", 3, 1, 63.64001366666667, 190.920041, 0.4959344089323702, 0.4969570394110899, 3, 1, 36.95112266666666, 110.853368, 0.1740115908809209, 0.1775751473829744, 2351.569951505343, 21164.12956354809, 0.08629833547091093, 0.08824721951643093),
  @("Inflammatory-Mac", "Tnf", "Tnfrsf1a", "FAPs", 3, 1, 63.64001366666667, 190.920041, 0.4959344089323702, 0.4969570394110899, 3, 1, 51.86084433333334, 155.582533, 0.2442250025331967, 0.2492264486514428, 3300.424842138206, 29703.82357924385, 0.1211195822778075, 0.123854838064761),
  @("Inflammatory-Mac", "Tnf", "Tnfrsf1a", "Inflammatory-Mac", 3, 1, 63.64001366666667, 190.920041, 0.4959344089323702, 0.4969570394110899, 3, 1, 64.73785366666665, 194.213561, 0.3048658902295037, 0.3111091917238571, 4119.917892097333, 37079.261028876, 0.1511934850746098, 0.1546079028526652),
  @("Inflammatory-Mac", "Tnf", "Tnfrsf1a", "Resolving-Mac", 3, 1, 63.64001366666667, 190.920041, 0.4959344089323702, 0.4969570394110899, 2, 1, 12.7841595, 25.568319, 0.06020363583370166, 0.04095769119761797, 813.5840852968466, 4881.50451178108, 0.02985705455276649, 0.02035421295868189),
  @("Inflammatory-Mac", "Tnf", "Tnfrsf1a", "MuSCs", 3, 1, 63.64001366666667, 190.920041, 0.4959344089323702, 0.4969570394110899, 3, 1, 46.01464833333333, 138.043945, 0.216693880522677, 0.2211315210441077, 2928.372848800194, 26355.35563920175, 0.1074659515562755, 0.1098928660185509),
  @("Resolving-Mac", "Tnf", "Tnfrsf1a", "ECs", 1, 0.5, 0.792186, 1.584372, 0.006173353414603005, 0.004124055360148531, 3, 1, 36.95112266666666, 110.853368, 0.1740115908809209, 0.1775751473829744, 29.27216206081599, 175.632972364896, 0.001074235048745234, 0.000732329738393921),
  @("Resolving-Mac", "Tnf", "Tnfrsf1a", "FAPs", 1, 0.5, 0.792186, 1.584372, 0.006173353414603005, 0.004124055360148531, 3, 1, 51.86084433333334, 155.582533, 0.2442250025331967, 0.2492264486514428, 41.08343482904601, 246.500608974276, 0.001507687253319737, 0.001027823671451765),
  @("Resolving-Mac", "Tnf", "Tnfrsf1a", "Inflammatory-Mac", 1, 0.5, 0.792186, 1.584372, 0.006173353414603005, 0.004124055360148531, 3, 1, 64.73785366666665, 194.213561, 0.3048658902295037, 0.3111091917238571, 51.28442134478198, 307.7065280686919, 0.001882044884444291, 0.00128303152972025),
  @("Resolving-Mac", "Tnf", "Tnfrsf1a", "Resolving-Mac", 1, 0.5, 0.792186, 1.584372, 0.006173353414603005, 0.004124055360148531, 2, 1, 12.7841595, 25.568319, 0.06020363583370166, 0.04095769119761797, 10.127432177667, 40.509728710668, 0.000371658320845498, 0.0001689117859228447),
  @("Resolving-Mac", "Tnf", "Tnfrsf1a", "MuSCs", 1, 0.5, 0.792186, 1.584372, 0.006173353414603005, 0.004124055360148531, 3, 1, 46.01464833333333, 138.043945, 0.216693880522677, 0.2211315210441077, 36.45216020459, 218.71296122754, 0.001337727907248244, 0.0009119586346597501),
  @("MuSCs", "Tnf", "Tnfrsf1a", "ECs", 3, 1, 57.101078, 171.303234, 0.444977738623016, 0.4458952950370734, 3, 1, 36.95112266666666, 110.853368, 0.1740115908809209, 0.1775751473829744, 2109.948937576901, 18989.54043819211, 0.07743128420438562, 0.07917992273358318),
  @("MuSCs", "Tnf", "Tnfrsf1a", "FAPs", 3, 1, 57.101078, 171.303234, 0.444977738623016, 0.4458952950370734, 3, 1, 51.86084433333334, 155.582533, 0.2442250025331967, 0.2492264486514428, 2961.310117423525, 26651.79105681173, 0.1086746893424222, 0.1111289008524771),
  @("MuSCs", "Tnf", "Tnfrsf1a", "Inflammatory-Mac", 3, 1, 57.101078, 171.303234, 0.444977738623016, 0.4458952950370734, 3, 1, 64.73785366666665, 194.213561, 0.3048658902295037, 0.3111091917238571, 3696.601231772919, 33269.41108595627, 0.1356585344176172, 0.1387221248324547),
  @("MuSCs", "Tnf", "Tnfrsf1a", "Resolving-Mac", 3, 1, 57.101078, 171.303234, 0.444977738623016, 0.4458952950370734, 2, 1, 12.7841595, 25.568319, 0.06020363583370166, 0.04095769119761797, 729.9892887739411, 4379.935732643647, 0.02678927773016414, 0.01826284180059921),
  @("MuSCs", "Tnf", "Tnfrsf1a", "MuSCs", 3, 1, 57.101078, 171.303234, 0.444977738623016, 0.4458952950370734, 3, 1, 46.01464833333333, 138.043945, 0.216693880522677, 0.2211315210441077, 2627.486023624237, 23647.37421261813, 0.09642395292842684, 0.09860150481795923)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $rowNum = $i + 2
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
  }
}
